$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9833462238311768
$ws.Range("B1").Value = 1.31148886680603
$ws.Range("C1").Value = 2.134623289108276
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 2.086764574050903
